# Efna3-Ephb1.xlsx update: refresh TPM-derived NATMI edge stats and add the two
# additional sending/target cluster combinations (MuSCs <-> ECs) that the new run produced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (Efna3/Ephb1)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna3"
$ws.Range("C2").Value = "Ephb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.046374
$ws.Range("H2").Value = 0.139122
$ws.Range("I2").Value = 0.6592866045237633
$ws.Range("J2").Value = 0.6592866045237632
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.475208
$ws.Range("N2").Value = 4.425624
$ws.Range("O2").Value = 0.8210007041987012
$ws.Range("P2").Value = 0.8210007041987013
$ws.Range("Q2").Value = 0.068411295792
$ws.Range("R2").Value = 0.6157016621279999
$ws.Range("S2").Value = 0.5412747665827803
$ws.Range("T2").Value = 0.5412747665827803

# Row 3: ECs -> MuSCs (Efna3/Ephb1)
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna3"
$ws.Range("C3").Value = "Ephb1"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.046374
$ws.Range("H3").Value = 0.139122
$ws.Range("I3").Value = 0.6592866045237633
$ws.Range("J3").Value = 0.6592866045237632
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3216333333333333
$ws.Range("N3").Value = 0.9649
$ws.Range("O3").Value = 0.1789992958012987
$ws.Range("P3").Value = 0.1789992958012987
$ws.Range("Q3").Value = 0.0149154242
$ws.Range("R3").Value = 0.1342388178
$ws.Range("S3").Value = 0.1180118379409829
$ws.Range("T3").Value = 0.1180118379409829

# Row 4: MuSCs -> ECs (Efna3/Ephb1)
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Efna3"
$ws.Range("C4").Value = "Ephb1"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.02396566666666667
$ws.Range("H4").Value = 0.071897
$ws.Range("I4").Value = 0.3407133954762367
$ws.Range("J4").Value = 0.3407133954762367
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.475208
$ws.Range("N4").Value = 4.425624
$ws.Range("O4").Value = 0.8210007041987012
$ws.Range("P4").Value = 0.8210007041987013
$ws.Range("Q4").Value = 0.035354343192
$ws.Range("R4").Value = 0.318189088728
$ws.Range("S4").Value = 0.279725937615921
$ws.Range("T4").Value = 0.279725937615921

# Row 5: MuSCs -> MuSCs (Efna3/Ephb1)
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Efna3"
$ws.Range("C5").Value = "Ephb1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.02396566666666667
$ws.Range("H5").Value = 0.071897
$ws.Range("I5").Value = 0.3407133954762367
$ws.Range("J5").Value = 0.3407133954762367
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3216333333333333
$ws.Range("N5").Value = 0.9649
$ws.Range("O5").Value = 0.1789992958012987
$ws.Range("P5").Value = 0.1789992958012987
$ws.Range("Q5").Value = 0.007708157255555555
$ws.Range("R5").Value = 0.0693734153
$ws.Range("S5").Value = 0.06098745786031576
$ws.Range("T5").Value = 0.06098745786031578
